{"js": "const replacements = [\n  [\"2026-02-21 Saturday\", \"2026-02-22 Sunday\"],\n  [\"495\u00f77=\", \"810\u00f79=\"],\n  [\"368\u00f78=\", \"575\u00f78=\"],\n  [\"775\u00f77=\", \"410\u00f79=\"],\n  [\"807\u00f77=\", \"158\u00f72=\"],\n  [\"781\u00f78=\", \"514\u00f77=\"],\n  [\"780\u00f78=\", \"378\u00f77=\"],\n  [\"673\u00f74=\", \"410\u00f78=\"],\n  [\"512\u00f79=\", \"274\u00f78=\"],\n  [\"675\u00f72=\", \"777\u00f76=\"],\n  [\"759\u00f72=\", \"156\u00f75=\"],\n  [\"211\u00f79=\", \"650\u00f76=\"],\n  [\"502\u00f78=\", \"997\u00f76=\"],\n  [\"821\u00f78=\", \"607\u00f75=\"],\n  [\"172\u00f72=\", \"319\u00f73=\"],\n  [\"784\u00f77=\", \"531\u00f79=\"],\n  [\"251\u00f75=\", \"377\u00f75=\"],\n  [\"811\u00f74=\", \"585\u00f79=\"],\n  [\"953\u00f79=\", \"313\u00f78=\"],\n  [\"623\u00f77=\", \"650\u00f77=\"],\n  [\"350\u00f77=\", \"830\u00f74=\"],\n  [\"812\u00f77=\", \"136\u00f74=\"],\n  [\"245\u00f74=\", \"246\u00f76=\"],\n  [\"703\u00f77=\", \"965\u00f79=\"],\n  [\"622\u00f75=\", \"923\u00f77=\"],\n  [\"480\u00f76=\", \"536\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-02-21 Saturday\", \"2026-02-22 Sunday\"),\n    @(\"495\u00f77=\", \"810\u00f79=\"),\n    @(\"368\u00f78=\", \"575\u00f78=\"),\n    @(\"775\u00f77=\", \"410\u00f79=\"),\n    @(\"807\u00f77=\", \"158\u00f72=\"),\n    @(\"781\u00f78=\", \"514\u00f77=\"),\n    @(\"780\u00f78=\", \"378\u00f77=\"),\n    @(\"673\u00f74=\", \"410\u00f78=\"),\n    @(\"512\u00f79=\", \"274\u00f78=\"),\n    @(\"675\u00f72=\", \"777\u00f76=\"),\n    @(\"759\u00f72=\", \"156\u00f75=\"),\n    @(\"211\u00f79=\", \"650\u00f76=\"),\n    @(\"502\u00f78=\", \"997\u00f76=\"),\n    @(\"821\u00f78=\", \"607\u00f75=\"),\n    @(\"172\u00f72=\", \"319\u00f73=\"),\n    @(\"784\u00f77=\", \"531\u00f79=\"),\n    @(\"251\u00f75=\", \"377\u00f75=\"),\n    @(\"811\u00f74=\", \"585\u00f79=\"),\n    @(\"953\u00f79=\", \"313\u00f78=\"),\n    @(\"623\u00f77=\", \"650\u00f77=\"),\n    @(\"350\u00f77=\", \"830\u00f74=\"),\n    @(\"812\u00f77=\", \"136\u00f74=\"),\n    @(\"245\u00f74=\", \"246\u00f76=\"),\n    @(\"703\u00f77=\", \"965\u00f79=\"),\n    @(\"622\u00f75=\", \"923\u00f77=\"),\n    @(\"480\u00f76=\", \"536\u00f76=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Not found: $oldText\"\n    }\n}\n"}
